$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-12 -> 2023-09-13) for every data row (rows 2 through 244).
$ws.Range("C2:C244").Value = 45182
